# MP119_Transform.xlsx refresh: re-run of the MATLAB export overwrote the
# X/Y/Z rotation values on Sheet1 and the workbook was saved with only the
# single data sheet (the two blank placeholder sheets are dropped).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.94726580371496549
$ws.Range("C2").Value = 0.2803986427237693
$ws.Range("D2").Value = -0.15512607218294502

$ws.Range("B3").Value = 0.29403616144410782
$ws.Range("C3").Value = -0.56810037548546366
$ws.Range("D3").Value = 0.76863821082254935

$ws.Range("B4").Value = 0.12739793120555307
$ws.Range("C4").Value = -0.77371736734543006
$ws.Range("D4").Value = -0.62058939935564639

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null
$excel.DisplayAlerts = $true

$ws.Range("A1:D4").Select() | Out-Null
